$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Value = 2.4
$ws.Range("Q2").Value = 1.66
$ws.Range("R2").Value = 1.55
$ws.Range("T2").Value = 1.74
$ws.Range("U2").Value = 2.26
$ws.Range("AL2").Value = 28
$ws.Range("AN2").Value = 7.6
$ws.Range("Q3").Value = 2.68
$ws.Range("V3").Value = 1.24
$ws.Range("G4").Value = 2.52
$ws.Range("W4").Value = 1.65
$ws.Range("P6").Value = 2.14
$ws.Range("Q6").Value = 1.66
$ws.Range("N7").Value = 1.88
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.1
$ws.Range("S7").Value = 2.12
$ws.Range("G8").Value = 1.65
$ws.Range("W8").Value = 2.52
$ws.Range("N9").Value = 2.96
$ws.Range("O9").Value = 1.49
$ws.Range("P9").Value = 1.65
$ws.Range("U9").Value = 1.68
$ws.Range("AA9").Value = 270
$ws.Range("AM9").Value = 280
$ws.Range("AO9").Value = 270
$ws.Range("F10").Value = 1.13
$ws.Range("J10").Value = 11.5
$ws.Range("K10").Value = 12.5
$ws.Range("L10").Value = 1.23
$ws.Range("N10").Value = 7.2
$ws.Range("O10").Value = 1.13
$ws.Range("P10").Value = 3.05
$ws.Range("Q10").Value = 1.37
$ws.Range("R10").Value = 1.82
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 2.68
$ws.Range("U10").Value = 1.46
$ws.Range("W10").Value = 8.2
$ws.Range("X10").Value = 46
$ws.Range("Y10").Value = 990
$ws.Range("AB10").Value = 13
$ws.Range("AC10").Value = 990
$ws.Range("AD10").Value = 990
$ws.Range("AG10").Value = 990
$ws.Range("AH10").Value = 990
$ws.Range("AJ10").Value = 7.6
$ws.Range("AL10").Value = 85
